$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "ДБТ София-Възраждане"
$ws.Range("B2").Value = "ДБТ София-Люлин"
$ws.Range("B3").Value = "ДБТ София-Сердика"
$ws.Range("B4").Value = "ДБТ София-Изток"
$ws.Range("B5").Value = "ДБТ Айтос"
$ws.Range("B6").Value = "ДБТ Бургас"
$ws.Range("B7").Value = "ДБТ Елхово"
$ws.Range("B8").Value = "ДБТ Карнобат"
$ws.Range("B9").Value = "ДБТ Поморие"
$ws.Range("B10").Value = "ДБТ Нова Загора"
$ws.Range("B11").Value = "ДБТ Сливен"
$ws.Range("B12").Value = "ДБТ Ямбол"
$ws.Range("B13").Value = "ДБТ Котел"
$ws.Range("B14").Value = "ДБТ Руен"
$ws.Range("B15").Value = "ДБТ Созопол"
$ws.Range("B16").Value = "ДБТ Варна"
$ws.Range("B17").Value = "ДБТ Вълчи Дол"
$ws.Range("B18").Value = "ДБТ Генерал Тошево"
$ws.Range("B19").Value = "ДБТ Долни Чифлик"
$ws.Range("B20").Value = "ДБТ Каварна"
$ws.Range("B21").Value = "ДБТ Каолиново"
$ws.Range("B22").Value = "ДБТ Нови Пазар "
$ws.Range("B23").Value = "ДБТ Велики Преслав"
$ws.Range("B24").Value = "ДБТ Провадия"
$ws.Range("B25").Value = "ДБТ Тервел"
$ws.Range("B26").Value = "ДБТ Добрич"
$ws.Range("B27").Value = "ДБТ Шумен"
$ws.Range("B28").Value = "ДБТ Велико Търново"
$ws.Range("B29").Value = "ДБТ Габрово"
$ws.Range("B30").Value = "ДБТ Горна Оряховица"
$ws.Range("B31").Value = "ДБТ Долна Митрополия"
$ws.Range("B32").Value = "ДБТ Левски"
$ws.Range("B33").Value = "ДБТ Ловеч"
$ws.Range("B34").Value = "ДБТ Луковит"
$ws.Range("B35").Value = "ДБТ Никопол"
$ws.Range("B36").Value = "ДБТ Павликени"
$ws.Range("B37").Value = "ДБТ Плевен"
$ws.Range("B38").Value = "ДБТ Свищов"
$ws.Range("B39").Value = "ДБТ Тетевен"
$ws.Range("B40").Value = "ДБТ Троян"
$ws.Range("B41").Value = "ДБТ Червен Бряг"
$ws.Range("B42").Value = "ДБТ Белоградчик"
$ws.Range("B43").Value = "ДБТ Берковица"
$ws.Range("B44").Value = "ДБТ Бяла Слатина"
$ws.Range("B45").Value = "ДБТ Видин"
$ws.Range("B46").Value = "ДБТ Враца"
$ws.Range("B47").Value = "ДБТ Козлодуй"
$ws.Range("B48").Value = "ДБТ Кула"
$ws.Range("B49").Value = "ДБТ Лом"
$ws.Range("B50").Value = "ДБТ Мездра"
$ws.Range("B51").Value = "ДБТ Монтана"
$ws.Range("B52").Value = "ДБТ Оряхово"
$ws.Range("B53").Value = "ДБТ Асеновград"
$ws.Range("B54").Value = "ДБТ Велинград"
$ws.Range("B55").Value = "ДБТ Девин"
$ws.Range("B56").Value = "ДБТ Златоград"
$ws.Range("B57").Value = "ДБТ Карлово"
$ws.Range("B58").Value = "ДБТ Мадан"
$ws.Range("B59").Value = "ДБТ Пловдив-Марица"
$ws.Range("B60").Value = "ДБТ Пазарджик"
$ws.Range("B61").Value = "ДБТ Панагюрище"
$ws.Range("B62").Value = "ДБТ Пещера"
$ws.Range("B63").Value = "ДБТ Пловдив"
$ws.Range("B64").Value = "ДБТ Първомай"
$ws.Range("B65").Value = "ДБТ Пловдив-Родопи"
$ws.Range("B66").Value = "ДБТ Септември"
$ws.Range("B67").Value = "ДБТ Смолян"
$ws.Range("B68").Value = "ДБТ Бяла"
$ws.Range("B69").Value = "ДБТ Ветово"
$ws.Range("B70").Value = "ДБТ Дулово"
$ws.Range("B71").Value = "ДБТ Исперих"
$ws.Range("B72").Value = "ДБТ Кубрат"
$ws.Range("B73").Value = "ДБТ Омуртаг"
$ws.Range("B74").Value = "ДБТ Попово"
$ws.Range("B75").Value = "ДБТ Разград"
$ws.Range("B76").Value = "ДБТ Сливница"
$ws.Range("B77").Value = "ДБТ Тутракан"
$ws.Range("B78").Value = "ДБТ Търговище"
$ws.Range("B79").Value = "ДБТ Русе"
$ws.Range("B80").Value = "ДБТ Благоевград"
$ws.Range("B81").Value = "ДБТ Ботевград"
$ws.Range("B82").Value = "ДБТ Гоце Делчев"
$ws.Range("B83").Value = "ДБТ Ихтиман"
$ws.Range("B84").Value = "ДБТ Кюстендил"
$ws.Range("B85").Value = "ДБТ Перник"
$ws.Range("B86").Value = "ДБТ Петрич"
$ws.Range("B87").Value = "ДБТ Разлог"
$ws.Range("B88").Value = "ДБТ Самоков"
$ws.Range("B89").Value = "ДБТ Сандански"
$ws.Range("B90").Value = "ДБТ Своге"
$ws.Range("B91").Value = "ДБТ Сливница"
$ws.Range("B92").Value = "ДБТ Пирдоп"
$ws.Range("B93").Value = "ДБТ Дупница"
$ws.Range("B94").Value = "ДБТ Ардино"
$ws.Range("B95").Value = "ДБТ Димитровград"
$ws.Range("B96").Value = "ДБТ Казанлък"
$ws.Range("B97").Value = "ДБТ Кирково"
$ws.Range("B98").Value = "ДБТ Крумовград"
$ws.Range("B99").Value = "ДБТ Кърджали"
$ws.Range("B100").Value = "ДБТ Момчилград"
$ws.Range("B101").Value = "ДБТ Раднево"
$ws.Range("B102").Value = "ДБТ Свиленград"
$ws.Range("B103").Value = "ДБТ Стара Загора"
$ws.Range("B104").Value = "ДБТ Харманли"
$ws.Range("B105").Value = "ДБТ Хасково"
$ws.Range("B106").Value = "ДБТ Чирпан"

$ws.Range("B80").Select()

